# Update Sheets via scheduled runner: refresh market-board derived profit values
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 271.63635
$ws.Range("J28").Value = 304.57144
$ws.Range("L28").Value = 304.57144
$ws.Range("N28").Value = -1274.57144
$ws.Range("H62").Value = 2315.8572
$ws.Range("I62").Value = 1950
$ws.Range("J62").Value = 2462.2
$ws.Range("K62").Value = 1950
$ws.Range("L62").Value = 2462.2
$ws.Range("M62").Value = -1326
$ws.Range("N62").Value = -3710.2
$ws.Range("H65").Value = 2315.8572
$ws.Range("I65").Value = 1950
$ws.Range("J65").Value = 2462.2
$ws.Range("K65").Value = 9750
$ws.Range("L65").Value = 12311
$ws.Range("M65").Value = -6630
$ws.Range("N65").Value = -18551
$ws.Range("H106").Value = 1295
$ws.Range("I106").Value = 1206.875
$ws.Range("K106").Value = 1206.875
$ws.Range("M106").Value = -575.875
$ws.Range("H124").Value = 69000
$ws.Range("J124").Value = 69000
$ws.Range("L124").Value = 69000
$ws.Range("N124").Value = -78820
$ws.Range("H129").Value = 996.6
$ws.Range("I129").Value = 749.625
$ws.Range("K129").Value = 2248.875
$ws.Range("M129").Value = 2751.125
$ws.Range("H137").Value = 1230.7646
$ws.Range("I137").Value = 1026.1351
$ws.Range("J137").Value = 1771.5714
$ws.Range("K137").Value = 3078.4053
$ws.Range("L137").Value = 5314.7142
$ws.Range("M137").Value = -528.4052999999999
$ws.Range("N137").Value = -10414.7142

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2188.577
$ws.Range("I61").Value = 1914.4286
$ws.Range("J61").Value = 3340
$ws.Range("K61").Value = 1914.4286
$ws.Range("L61").Value = 3340
$ws.Range("M61").Value = -1702.4286
$ws.Range("N61").Value = -3764
$ws.Range("H88").Value = 2503.3333
$ws.Range("I88").Value = 1506
$ws.Range("J88").Value = 2702.8
$ws.Range("K88").Value = 1506
$ws.Range("L88").Value = 2702.8
$ws.Range("M88").Value = -1100
$ws.Range("N88").Value = -3514.8
$ws.Range("H91").Value = 2503.3333
$ws.Range("I91").Value = 1506
$ws.Range("J91").Value = 2702.8
$ws.Range("K91").Value = 1506
$ws.Range("L91").Value = 2702.8
$ws.Range("M91").Value = -102
$ws.Range("N91").Value = -5510.8
$ws.Range("H97").Value = 1500
$ws.Range("I97").Value = 1500
$ws.Range("K97").Value = 1500
$ws.Range("M97").Value = -1004
$ws.Range("H134").Value = 50266.668
$ws.Range("J134").Value = 50266.668
$ws.Range("L134").Value = 50266.668
$ws.Range("N134").Value = -60406.668
$ws.Range("H136").Value = 2188.577
$ws.Range("I136").Value = 1914.4286
$ws.Range("J136").Value = 3340
$ws.Range("K136").Value = 5743.2858
$ws.Range("L136").Value = 10020
$ws.Range("M136").Value = -3193.2858
$ws.Range("N136").Value = -15120

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 98742.5
$ws.Range("J57").Value = 98742.5
$ws.Range("L57").Value = 98742.5
$ws.Range("N57").Value = -100182.5
$ws.Range("H86").Value = 149017.14
$ws.Range("I86").Value = 4435.3335
$ws.Range("J86").Value = 257453.5
$ws.Range("K86").Value = 4435.3335
$ws.Range("L86").Value = 257453.5
$ws.Range("M86").Value = -3312.3335
$ws.Range("N86").Value = -259699.5
$ws.Range("H89").Value = 149017.14
$ws.Range("I89").Value = 4435.3335
$ws.Range("J89").Value = 257453.5
$ws.Range("K89").Value = 22176.6675
$ws.Range("L89").Value = 1287267.5
$ws.Range("M89").Value = -16560.6675
$ws.Range("N89").Value = -1298499.5
$ws.Range("H105").Value = 3291.074
$ws.Range("I105").Value = 3028.842
$ws.Range("J105").Value = 3913.875
$ws.Range("K105").Value = 3028.842
$ws.Range("L105").Value = 3913.875
$ws.Range("M105").Value = -1281.842
$ws.Range("N105").Value = -7407.875
$ws.Range("H134").Value = 1909.8
$ws.Range("I134").Value = 1874.7307
$ws.Range("K134").Value = 5624.1921
$ws.Range("M134").Value = -3089.1921
$ws.Range("H136").Value = 98742.5
$ws.Range("J136").Value = 98742.5
$ws.Range("L136").Value = 98742.5
$ws.Range("N136").Value = -108942.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 89.55556
$ws.Range("I7").Value = 47.076923
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 47.076923
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = 65.92307700000001
$ws.Range("N7").Value = -426
$ws.Range("H31").Value = 1904.3148
$ws.Range("I31").Value = 1362.2195
$ws.Range("K31").Value = 1362.2195
$ws.Range("M31").Value = -1067.2195
$ws.Range("H34").Value = 1904.3148
$ws.Range("I34").Value = 1362.2195
$ws.Range("K34").Value = 1362.2195
$ws.Range("M34").Value = -1160.2195
$ws.Range("H58").Value = 1690.3572
$ws.Range("I58").Value = 1967.3636
$ws.Range("J58").Value = 674.6667
$ws.Range("K58").Value = 1967.3636
$ws.Range("L58").Value = 674.6667
$ws.Range("M58").Value = -1764.3636
$ws.Range("N58").Value = -1080.6667
$ws.Range("H99").Value = 2668.9
$ws.Range("I99").Value = 2811.5386
$ws.Range("K99").Value = 2811.5386
$ws.Range("M99").Value = -1313.5386
$ws.Range("H122").Value = 2430.3872
$ws.Range("I122").Value = 2201.4
$ws.Range("K122").Value = 6604.200000000001
$ws.Range("M122").Value = -4154.200000000001
$ws.Range("H126").Value = 2668.9
$ws.Range("I126").Value = 2811.5386
$ws.Range("K126").Value = 8434.6158
$ws.Range("M126").Value = -5964.6158
$ws.Range("H132").Value = 4333
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 12999
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -18059
$ws.Range("H134").Value = 2137.84
$ws.Range("I134").Value = 1945
$ws.Range("K134").Value = 5835
$ws.Range("M134").Value = -3300
$ws.Range("H136").Value = 1690.3572
$ws.Range("I136").Value = 1967.3636
$ws.Range("J136").Value = 674.6667
$ws.Range("K136").Value = 5902.0908
$ws.Range("L136").Value = 2024.0001
$ws.Range("M136").Value = -3352.0908
$ws.Range("N136").Value = -7124.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 614.0263
$ws.Range("I122").Value = 435.18182
$ws.Range("J122").Value = 686.8889
$ws.Range("K122").Value = 3916.63638
$ws.Range("L122").Value = 6182.0001
$ws.Range("M122").Value = -1466.63638
$ws.Range("N122").Value = -11082.0001
$ws.Range("H132").Value = 961.5714
$ws.Range("I132").Value = 1020.2
$ws.Range("K132").Value = 9181.800000000001
$ws.Range("M132").Value = -6651.800000000001
$ws.Range("H137").Value = 3030
$ws.Range("I137").Value = 3030
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 9090
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3990
$ws.Range("N137").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 102.30769
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 118.181816
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 118.181816
$ws.Range("M2").Value = 98
$ws.Range("N2").Value = -344.181816
$ws.Range("H80").Value = 3150
$ws.Range("I80").Value = 3062.5
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 3062.5
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -2064.5
$ws.Range("N80").Value = -5496
$ws.Range("H83").Value = 3150
$ws.Range("I83").Value = 3062.5
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 15312.5
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -10320.5
$ws.Range("N83").Value = -27484
$ws.Range("H122").Value = 1951.5454
$ws.Range("I122").Value = 1058.375
$ws.Range("J122").Value = 4333.3335
$ws.Range("K122").Value = 3175.125
$ws.Range("L122").Value = 13000.0005
$ws.Range("M122").Value = -725.125
$ws.Range("N122").Value = -17900.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 37506160
$ws.Range("I122").Value = 41672316
$ws.Range("J122").Value = 33340000
$ws.Range("K122").Value = 125016948
$ws.Range("L122").Value = 100020000
$ws.Range("M122").Value = -125014498
$ws.Range("N122").Value = -100024900
$ws.Range("H132").Value = 4108.729
$ws.Range("I132").Value = 3596.6216
$ws.Range("J132").Value = 4970
$ws.Range("K132").Value = 10789.8648
$ws.Range("L132").Value = 14910
$ws.Range("M132").Value = -8259.864799999999
$ws.Range("N132").Value = -19970

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10778211
$ws.Range("I122").Value = 13160058
$ws.Range("K122").Value = 39480174
$ws.Range("M122").Value = -39477724
$ws.Range("H132").Value = 1833.1666
$ws.Range("I132").Value = 1075.421
$ws.Range("J132").Value = 3142
$ws.Range("K132").Value = 3226.263
$ws.Range("L132").Value = 9426
$ws.Range("M132").Value = -696.2629999999999
$ws.Range("N132").Value = -14486
